$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; existing rows 31-46 shift down to 32-47.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with its data (weekly price update).
$ws.Cells.Item(31, 1).Value = 11
$ws.Cells.Item(31, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value = "Bíobío"
$ws.Cells.Item(31, 4).Value = 44489
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(31, 6).Value = 100112012
$ws.Cells.Item(31, 7).Value = "Espinaca"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 50
$ws.Cells.Item(31, 11).Value = 6000
$ws.Cells.Item(31, 12).Value = 6500
$ws.Cells.Item(31, 13).Value = 6300
$ws.Cells.Item(31, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(31, 15).Value = "Región Metropolitana"
$ws.Cells.Item(31, 16).Value = 630
$ws.Cells.Item(31, 17).Value = 10
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Match the date cell's number format to the rest of column D (dates above/below it).
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
